$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing rows down
$ws.Rows.Item(1).Insert()

# Populate the newly inserted row 1 with the new data point
# (leading apostrophe forces text, matching the date-as-text values used
# throughout column A)
$ws.Range("A1").Value = "'20181011"
$ws.Range("B1").Value = 958.0
